# Tableau_charges.xlsx - "Maj docs suivi de projet" (suivi perso / suivi taches)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- suivi taches: fill in hours consumed (column D) for rows 44-46 ---
# (F = C - D recalculates automatically, as do the SUM/AVERAGE rows above)
$ws.Range("D44").Value = 6
$ws.Range("D45").Value = 10
$ws.Range("D46").Value = 5

# --- suivi perso: move the current selection to D47 ---
$ws.Range("D47").Select()
